$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl_spp")

# The table was previously sorted by column A (spc); re-sort it by
# column B (spc_grp_orig) instead, keeping the header row in place.
$lo = $ws.ListObjects.Item("Table1")
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("B1:B122"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Hide the now-redundant spc_grp_alt (C) and helper orig1_same_as_alt1 (F)
# columns.
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(6).Hidden = $true

# Leave the selection on the (now hidden) helper column, matching the
# state Excel leaves behind after hiding the selected column.
$ws.Activate()
$ws.Columns.Item(6).Select()
